$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "2025-04-28 09:15:10"
$ws.Range("B8").Value = 204
